$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 397; existing rows 397-408 shift down to 400-411.
$ws.Rows("397:399").Insert()

# Populate the three newly inserted rows (397-399) with the new weekly data.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Categoria ID,
#          G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo,
#          L Precio maximo, M Precio promedio ponderado, N Unidad de comercializacion,
#          O Origen, P Precio $/Kg, Q Kg o Unidades, R Clasificacion

$ws.Range("A397:A399").Value = 10
$ws.Range("B397:B399").Value = "Vega Modelo de Temuco"
$ws.Range("C397:C399").Value = "La Araucanía"
$ws.Range("E397:E399").Value = 9
$ws.Range("F397:F399").Value = 100114001
$ws.Range("G397:G399").Value = "Papa"
$ws.Range("O397:O399").Value = "Provincia de Cautín"
$ws.Range("Q397:Q399").Value = 25
$ws.Range("R397:R399").Value = "Hortaliza"

$ws.Range("D397").Value = 44448
$ws.Range("H397").Value = "Asterix"
$ws.Range("I397").Value = "1a (guarda)"
$ws.Range("J397").Value = 850
$ws.Range("K397").Value = 8000
$ws.Range("L397").Value = 8000
$ws.Range("M397").Value = 8000
$ws.Range("N397").Value = "$/malla 25 kilos"
$ws.Range("P397").Value = 320

$ws.Range("D398").Value = 44448
$ws.Range("H398").Value = "Asterix"
$ws.Range("I398").Value = "1a (guarda)"
$ws.Range("J398").Value = 650
$ws.Range("K398").Value = 7000
$ws.Range("L398").Value = 7000
$ws.Range("M398").Value = 7000
$ws.Range("N398").Value = "$/saco 25 kilos"
$ws.Range("P398").Value = 280

$ws.Range("D399").Value = 44448
$ws.Range("H399").Value = "Rosara"
$ws.Range("I399").Value = "1a (guarda)"
$ws.Range("J399").Value = 550
$ws.Range("K399").Value = 7000
$ws.Range("L399").Value = 7000
$ws.Range("M399").Value = 7000
$ws.Range("N399").Value = "$/saco 25 kilos"
$ws.Range("P399").Value = 280
